# Scheduled runner update: refresh cached market-price-derived columns
# (currentAveragePrice / currentAveragePriceNQ / currentAveragePriceHQ /
# LevePriceNQ / LevePriceHQ / LeveProfitNQ / LeveProfitHQ) for a handful
# of leves across every crafting-job sheet.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H11").Value = 128.71428
$ws.Range("I11").Value = 128.71428
$ws.Range("K11").Value = 128.71428
$ws.Range("M11").Value = 11.28572

$ws.Range("H51").Value = 6445.05
$ws.Range("I51").Value = 4600.125
$ws.Range("J51").Value = 7675
$ws.Range("K51").Value = 4600.125
$ws.Range("L51").Value = 7675
$ws.Range("M51").Value = -4116.125
$ws.Range("N51").Value = -8643

$ws.Range("H70").Value = 2883.3572
$ws.Range("I70").Value = 3799.1667
$ws.Range("J70").Value = 1234.9
$ws.Range("K70").Value = 11397.5001
$ws.Range("L70").Value = 3704.7
$ws.Range("M70").Value = -11127.5001
$ws.Range("N70").Value = -4244.700000000001

$ws.Range("H73").Value = 2883.3572
$ws.Range("I73").Value = 3799.1667
$ws.Range("J73").Value = 1234.9
$ws.Range("K73").Value = 11397.5001
$ws.Range("L73").Value = 3704.7
$ws.Range("M73").Value = -10461.5001
$ws.Range("N73").Value = -5576.700000000001

$ws.Range("H116").Value = 40983.25
$ws.Range("I116").Value = 61625.11
$ws.Range("J116").Value = 3827.9
$ws.Range("K116").Value = 61625.11
$ws.Range("L116").Value = 3827.9
$ws.Range("M116").Value = -58183.11
$ws.Range("N116").Value = -10711.9

$ws.Range("H136").Value = 0
$ws.Range("J136").Value = 0
$ws.Range("L136").Value = 0
$ws.Range("N136").ClearContents()

$ws.Range("H139").Value = 30000
$ws.Range("I139").Value = 30000
$ws.Range("J139").Value = 0
$ws.Range("K139").Value = 30000
$ws.Range("L139").Value = 0
$ws.Range("M139").Value = -24860
$ws.Range("N139").ClearContents()

$ws.Range("H140").Value = 0
$ws.Range("J140").Value = 0
$ws.Range("L140").Value = 0
$ws.Range("N140").ClearContents()

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H9").Value = 58333.332
$ws.Range("I9").Value = 55000
$ws.Range("K9").Value = 55000
$ws.Range("M9").Value = -54830

$ws.Range("H20").Value = 58333.332
$ws.Range("I20").Value = 55000
$ws.Range("K20").Value = 55000
$ws.Range("M20").Value = -54730

$ws.Range("H63").Value = 11123.4
$ws.Range("I63").Value = 11123.4
$ws.Range("K63").Value = 11123.4
$ws.Range("M63").Value = -10437.4

$ws.Range("H66").Value = 11123.4
$ws.Range("I66").Value = 11123.4
$ws.Range("K66").Value = 55617
$ws.Range("M66").Value = -52185

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 41324.957
$ws.Range("I105").Value = 80580.836
$ws.Range("J105").Value = 2069.0833
$ws.Range("K105").Value = 80580.836
$ws.Range("L105").Value = 2069.0833
$ws.Range("M105").Value = -78833.836
$ws.Range("N105").Value = -5563.0833

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 872.125
$ws.Range("I16").Value = 802.2
$ws.Range("K16").Value = 802.2
$ws.Range("M16").Value = -515.2

$ws.Range("H86").Value = 2692.7104
$ws.Range("I86").Value = 1869.7084
$ws.Range("J86").Value = 4103.5713
$ws.Range("K86").Value = 1869.7084
$ws.Range("L86").Value = 4103.5713
$ws.Range("M86").Value = -746.7084
$ws.Range("N86").Value = -6349.5713

$ws.Range("H89").Value = 2692.7104
$ws.Range("I89").Value = 1869.7084
$ws.Range("J89").Value = 4103.5713
$ws.Range("K89").Value = 9348.541999999999
$ws.Range("L89").Value = 20517.8565
$ws.Range("M89").Value = -3732.541999999999
$ws.Range("N89").Value = -31749.8565

$ws.Range("H113").Value = 872.125
$ws.Range("I113").Value = 802.2
$ws.Range("K113").Value = 802.2
$ws.Range("M113").Value = 1367.8

$ws.Range("H122").Value = 2959.9
$ws.Range("J122").Value = 1800
$ws.Range("L122").Value = 5400
$ws.Range("N122").Value = -10300

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H68").Value = 1099.2069
$ws.Range("I68").Value = 633.69446
$ws.Range("J68").Value = 1860.9546
$ws.Range("K68").Value = 1901.08338
$ws.Range("L68").Value = 5582.8638
$ws.Range("M68").Value = -1090.08338
$ws.Range("N68").Value = -7204.8638

$ws.Range("H71").Value = 1099.2069
$ws.Range("I71").Value = 633.69446
$ws.Range("J71").Value = 1860.9546
$ws.Range("K71").Value = 5703.25014
$ws.Range("L71").Value = 16748.5914
$ws.Range("M71").Value = -1647.25014
$ws.Range("N71").Value = -24860.5914

$ws.Range("H92").Value = 492.75
$ws.Range("I92").Value = 324
$ws.Range("J92").Value = 999
$ws.Range("K92").Value = 972
$ws.Range("L92").Value = 2997
$ws.Range("M92").Value = 276
$ws.Range("N92").Value = -5493

$ws.Range("H113").Value = 620.4286
$ws.Range("I113").Value = 603
$ws.Range("J113").Value = 643.6667
$ws.Range("K113").Value = 1809
$ws.Range("L113").Value = 1931.0001
$ws.Range("M113").Value = 361
$ws.Range("N113").Value = -6271.0001

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H19").Value = 584.2857
$ws.Range("I19").Value = 584.2857
$ws.Range("J19").Value = 0
$ws.Range("K19").Value = 584.2857
$ws.Range("L19").Value = 0
$ws.Range("M19").Value = -296.2857
$ws.Range("N19").ClearContents()

$ws.Range("H113").Value = 29418984
$ws.Range("J113").Value = 1755.7778
$ws.Range("L113").Value = 1755.7778
$ws.Range("N113").Value = -6095.7778

$ws.Range("H122").Value = 1948.8
$ws.Range("I122").Value = 1800
$ws.Range("J122").Value = 1986
$ws.Range("K122").Value = 5400
$ws.Range("L122").Value = 5958
$ws.Range("M122").Value = -2950
$ws.Range("N122").Value = -10858

$ws.Range("H131").Value = 24884
$ws.Range("J131").Value = 24884
$ws.Range("L131").Value = 24884
$ws.Range("N131").Value = -34964

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 1053.75
$ws.Range("I61").Value = 997.5833
$ws.Range("J61").Value = 1222.25
$ws.Range("K61").Value = 997.5833
$ws.Range("L61").Value = 1222.25
$ws.Range("M61").Value = -795.5833
$ws.Range("N61").Value = -1626.25

$ws.Range("H113").Value = 1053.75
$ws.Range("I113").Value = 997.5833
$ws.Range("J113").Value = 1222.25
$ws.Range("K113").Value = 997.5833
$ws.Range("L113").Value = 1222.25
$ws.Range("M113").Value = 1172.4167
$ws.Range("N113").Value = -5562.25

$ws.Range("H139").Value = 0
$ws.Range("J139").Value = 0
$ws.Range("L139").Value = 0
$ws.Range("N139").ClearContents()

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 4400.154
$ws.Range("J81").Value = 5666.8887
$ws.Range("L81").Value = 11333.7774
$ws.Range("N81").Value = -13455.7774

$ws.Range("H84").Value = 4400.154
$ws.Range("J84").Value = 5666.8887
$ws.Range("L84").Value = 56668.887
$ws.Range("N84").Value = -67276.887

$ws.Range("H113").Value = 386.42105
$ws.Range("I113").Value = 452.36365
$ws.Range("J113").Value = 295.75
$ws.Range("K113").Value = 1357.09095
$ws.Range("L113").Value = 887.25
$ws.Range("M113").Value = 812.90905
$ws.Range("N113").Value = -5227.25

$ws.Range("H131").Value = 23900
$ws.Range("J131").Value = 23900
$ws.Range("L131").Value = 23900
$ws.Range("N131").Value = -33980
